$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.517
$ws.Range("C4").Value = -12.88
$ws.Range("A7").Value = -19.796
$ws.Range("C12").Value = -10.662
$ws.Range("A16").Value = -22.059
$ws.Range("C18").Value = -12.176
$ws.Range("C19").Value = -11.77
$ws.Range("C20").Value = -12.067
$ws.Range("A28").Value = -21.934
$ws.Range("A29").Value = -21.344
$ws.Range("C31").Value = -13.298
$ws.Range("A32").Value = -21.781
$ws.Range("A40").Value = -20.234
$ws.Range("C40").Value = -12.226
$ws.Range("C42").Value = -12.355
$ws.Range("C47").Value = -11.844
$ws.Range("C48").Value = -11.897
$ws.Range("A52").Value = -21.918
$ws.Range("A57").Value = -22.253
$ws.Range("C63").Value = -11.207
$ws.Range("C64").Value = -10.739
$ws.Range("A66").Value = -21.547
$ws.Range("C76").Value = -12.969
$ws.Range("C81").Value = -13.1
$ws.Range("C89").Value = -13.51
$ws.Range("C94").Value = -11.539
$ws.Range("A100").Value = -22.157
